$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 54, shifting existing rows 54-66 down to 55-67
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly record
$ws.Range("A54").Value = 5
$ws.Range("B54").Value = "Macroferia Regional de Talca"
$ws.Range("C54").Value = "Maule"
$ws.Range("D54").Value = 44522
$ws.Range("E54").Value = 7
$ws.Range("F54").Value = 100112022
$ws.Range("G54").Value = "Arveja Verde"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 500
$ws.Range("K54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("M54").Value = 15000
$ws.Range("N54").Value = "$/saco 25 kilos"
$ws.Range("O54").Value = "Región del Maule"
$ws.Range("P54").Value = 600
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"
